$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# New header cell J7: "Loading Details Name" - reuse existing header style from I7
$hdr = $ws.Cells.Item(7, 10)
$hdr.Value = "Loading Details Name"
$ws.Range("I7").Copy()
$hdr.PasteSpecial(-4122)

# New data cells J8:J14: "40V (A)" - reuse fill/border style from column A data cells,
# then add word-wrap to match the new style used in the diff
for ($r = 8; $r -le 14; $r++) {
  $cell = $ws.Cells.Item($r, 10)
  $cell.Value = "40V (A)"
  $ws.Range("A8").Copy()
  $cell.PasteSpecial(-4122)
  $cell.WrapText = $true
}

$ws.Application.CutCopyMode = $false

# Autofit the new column to size it to its content
$ws.Columns.Item(10).EntireColumn.AutoFit()

# Update the selection to match the new active cell / selected range
[void]$ws.Range("J9:J14").Select()
